$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.393.87"
$ws.Range("E2").Value = "  +6.68%  "
$ws.Range("D3").Value = "3.562.25"
$ws.Range("E3").Value = "  +3.89%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'419.30"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").Value = "'132.31"
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("D7").Value = "'0.661"
$ws.Range("E7").Value = "  +6.15%  "
$ws.Range("D8").Value = "3.551.69"
$ws.Range("E8").Value = "  +3.77%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.786"
$ws.Range("E10").Value = "  +8.36%  "
$ws.Range("E11").Value = "  +20.92%  "
$ws.Range("D12").Value = "'0.0000290"
$ws.Range("E12").Value = "  +34.05%  "
$ws.Range("D13").Value = "'43.62"
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").Value = "'10.15"
$ws.Range("E14").Value = "  +9.33%  "
$ws.Range("D15").Value = "4.117.57"
$ws.Range("E15").Value = "  +3.69%  "
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "'20.54"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "3.533.44"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("E19").Value = "  +5.02%  "
$ws.Range("D20").Value = "'12.77"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "66.212.19"
$ws.Range("E21").Value = "  +6.25%  "
$ws.Range("D22").Value = "'450.52"
$ws.Range("E22").Value = "  -3.64%  "
$ws.Range("D23").Value = "'90.48"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D25").Value = "'13.21"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").Value = "'3.39"
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("D27").Value = "'10.09"
$ws.Range("E27").Value = "  -4.13%  "
$ws.Range("D28").Value = "'34.45"
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("D29").Value = "'4.83"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "'2.81"
$ws.Range("E30").Value = "  +6.36%  "
$ws.Range("D31").Value = "'12.51"
$ws.Range("E31").Value = "  +4.77%  "
$ws.Range("E32").Value = "  +6.22%  "
$ws.Range("D33").Value = "'7.33"
$ws.Range("E33").Value = "  -4.41%  "
$ws.Range("D34").Value = "'0.163"
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'39.19"
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("D38").Value = "'0.0508"
$ws.Range("E38").Value = "  +4.21%  "
$ws.Range("D39").Value = "0.0₃0741"
$ws.Range("E39").Value = "  +39.68%  "
$ws.Range("E40").Value = "  +10.91%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.07"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'0.996"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  +4.09%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'4.48"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'147.45"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").Value = "'3.28"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Value = "'0.312"
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("E49").Value = "  -4.78%  "
$ws.Range("E50").Value = "  +6.75%  "
$ws.Range("D51").Value = "'15.76"
$ws.Range("E51").Value = "  -3.94%  "
